$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.988.20"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "3.147.52"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.79"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.149.01"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.09"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").Value = "3.664.27"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "64.928.06"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.13"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("D18").Value = "3.140.84"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "502.30"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("E22").Value = "  -3.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.10"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.72"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.18"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.05"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.60%  "
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.80"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.45"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.44"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.62%  "
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.47"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.97"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0886"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "474.84"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("E39").Value = "  -2.53%  "
$ws.Range("E40").Value = "  -3.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.72"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("D42").Value = "3.000.28"
$ws.Range("E42").Value = "  -3.90%  "
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.281"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.41%  "
$ws.Range("E45").Value = "  -2.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.24"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.54%  "
$ws.Range("D47").Value = "0.0₃0596"
$ws.Range("E47").Value = "  +2.92%  "
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("E50").Value = "  -2.76%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.91"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.01%  "
